# Applies the "added envdetails to plot panels loadtab" edit:
#  - Sheet "Example Test-1": extend the blank-but-styled data cells from the
#    Load-5 column (F) out through Load-9/Load-Max (G:L) for the measurement
#    rows (6-14), and add the same blank-but-styled cells across B:L for the
#    unit rows (15-22). No values are entered on this sheet - only formatting
#    so the cells behave like their already-populated neighbours.
#  - Sheet "Example Test-2": same formatting extension (this time starting
#    from column H, since F:G were already present) plus a handful of actual
#    "Load-Max" (column L) and "T" row (21) sample values.
#  - Update the remembered selection on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "Example Test-1" ---------------------------------------------
# Use an existing blank formatted cell (F7, style shared by every blank
# measurement cell on the sheet) as the format donor.
$fmtSrc1 = $ws1.Range("F7")
$fmtSrc1.Copy()
$ws1.Range("G6:L14").PasteSpecial(-4122)

$fmtSrc1.Copy()
$ws1.Range("B15:L22").PasteSpecial(-4122)

# --- Sheet "Example Test-2" ---------------------------------------------
$fmtSrc2 = $ws2.Range("F7")
$fmtSrc2.Copy()
$ws2.Range("H6:L14").PasteSpecial(-4122)

$fmtSrc2.Copy()
$ws2.Range("B15:L22").PasteSpecial(-4122)

# Load-Max (column L) sample readings added alongside the rest of the row.
$ws2.Range("L6").Value = 160
$ws2.Range("L9").Value = 3.2
$ws2.Range("L10").Value = 140
$ws2.Range("L11").Value = 95
$ws2.Range("L14").Value = 40

# New "T" (temperature) row sample values across Load-Rest..Load-5, plus Load-Max.
$ws2.Range("B21").Value = 37
$ws2.Range("C21").Value = 37.1
$ws2.Range("D21").Value = 37.2
$ws2.Range("E21").Value = 37.3
$ws2.Range("F21").Value = 37.4
$ws2.Range("G21").Value = 37.5
$ws2.Range("L21").Value = 38

# --- Remembered selections -----------------------------------------------
$ws2.Range("G13").Select()
$ws1.Activate()
$ws1.Range("H11").Select()
